$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "01‏/05‏/2025 02:11:35 م"
$ws.Range("B4").Value = "NRC"
$ws.Range("C4").Value = "C2"
$ws.Range("D4").Value = "الرحلة 2"
$ws.Range("E4").Value = "بير 19"
$ws.Range("F4").Value = "يامن "
$ws.Range("G4").Value = "'2323"
$ws.Range("H4").Value = "واو"

# G3 already holds a plain (non quote-prefixed) text value; copy its
# formatting onto G4 so the new cell keeps the same style as its
# neighbours instead of picking up a distinct "quote prefix" style.
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)
